$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (B13/C13 = "1285870 - Marcos Villela Barcza", no A13 label) is
# removed entirely; every row below it shifts up by one. Row heights and the
# column-A labels all line up correctly after a plain row delete, so no
# further per-row restyling is required.
$ws.Rows.Item(13).Delete()

# After the shift a handful of rows end up with content that doesn't match
# the label that landed on them (this mirrors the source data exactly) -
# patch those B/C pairs explicitly.

# Row 10 "Objetivos:" now shows the Docentes responsaveis name instead of the
# old mission-statement paragraph.
$ws.Range("B10").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C10").Value = "1285870 - Marcos Villela Barcza"

# Row 13 "Programa resumido:" now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 "Programa:" now shows the activation date instead of the long
# syllabus paragraph.
$ws.Range("B15").Value = "15/07/2015"
$ws.Range("C15").Value = "15/07/2015"

# Row 18 "Metodo:" now shows the Docentes responsaveis name.
$ws.Range("B18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C18").Value = "1285870 - Marcos Villela Barcza"

# Row 19 "Criterio:" now shows the teaching-method description.
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários"
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários"

# Row 20 "Norma de recuperação:" now shows the evaluation-criteria text.
$ws.Range("B20").Value = "Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula."

# Row 21 "Bibliografia:" now shows the recovery-exam rule instead of the book
# list.
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
